$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 data: B2 -> "Admin", C2 -> "admin123" (as text)
$ws.Range("B2").Value = "Admin"
$ws.Range("C2").Value = "admin123"

# Update selection to B3
[void]$ws.Range("B3").Select()

# Set page orientation (portrait) to add pageSetup element
$ws.PageSetup.Orientation = 1
